# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the existing header cells (e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2 through 46)
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}
